# [MOSIP-43633] Added packetCreatedOn to identity schema for biometric issue fix
#
# The identity schema JSON lives in cell F2 ("schema_json" column) of Sheet1.
# We insert a new "packetCreatedOn" property (a clone of the existing
# "typeOfDeath" property) right after "typeOfDeath" inside
# properties.identity.properties.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("F2")
$json = $cell.Value2

$oldBlock = @"
                "typeOfDeath": {
                    "bioAttributes": [
                    ],
                    "fieldCategory": "pvt",
                    "format": "none",
                    "type": "string",
                    "fieldType": "default"
                }
"@

$newBlock = @"
                "typeOfDeath": {
                    "bioAttributes": [
                    ],
                    "fieldCategory": "pvt",
                    "format": "none",
                    "type": "string",
                    "fieldType": "default"
                },
`t`t`t`t"packetCreatedOn": {
                    "bioAttributes": [
                    ],
                    "fieldCategory": "pvt",
                    "format": "none",
                    "type": "string",
                    "fieldType": "default"
                }
"@

if ($json.IndexOf($oldBlock) -lt 0) {
    throw "Could not locate the 'typeOfDeath' field block in the schema_json cell; aborting edit."
}

$updated = $json.Replace($oldBlock, $newBlock)

$cell.Value2 = $updated

# Writing a longer value into this word-wrapped cell makes the host
# recompute the row's autofit height; restore the original (Excel-capped)
# row height so row 2 keeps looking the way it did before the edit.
$ws.Rows.Item(2).RowHeight = 409.5
